$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H (rows 2-27) currently stores a boolean FALSE rendered through a
# custom "TRUE"/"FALSE" number format. Replace each of those boolean cells
# with the literal text value "False" (a real text string, not a boolean),
# formatted with the standard Text number format.
$rng = $ws.Range("H2:H27")

# A leading apostrophe forces Excel to store the value as literal text
# instead of auto-coercing the word "False" back into a boolean.
$rng.Value = "'False"

# Drop back to the default style, then apply the Text number format so the
# cells end up with a clean (non quote-prefixed) text format, matching a
# plain "store as text" number format rather than the old boolean format.
$rng.ClearFormats()
$rng.NumberFormat = "@"

# Match the saved selection/active cell.
$ws.Range("H3:H27").Select()
